# Update "想去人数" (column F, rows 2-9) values on the "展览" and "全部类型"
# worksheets to reflect newly generated output figures.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$updates = @{
    2 = 2205
    3 = 1670
    4 = 325
    5 = 1071
    6 = 683
    7 = 35
    8 = 5765
    9 = 85
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
